$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newPattern = "hhhhhhhh-hhhh-hhhh-hhhh-hhhhhhhhhhhh"

$ws.Range("G19").Value = $newPattern
$ws.Range("G24").Value = $newPattern
$ws.Range("G29").Value = $newPattern
$ws.Range("G34").Value = $newPattern
